$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Content") - this shifts C:F left to B:E
$ws.Range("B1").EntireColumn.Delete()

# Update header and value in column A
$ws.Range("A1").Value = "Student ID"
$ws.Range("A2").Value = "Gg"
